# Updated capital structure database
# - Row 2 (stub row "1"): refreshed ratio columns, dropped the two
#   historical-growth columns (D/E) and the buybacks_cash_returned column (T),
#   added roe (W) / roe_cost_equity (Y).
# - Row 3: company renamed from "Standard Chartered Bank Zambia Plc (LUSE:SCZ)"
#   to "Cavmont Capital Holdings Zambia Plc (LUSE:CCHZ)" with the same refresh.
# - Row 4 (Zambia National Commercial Bank Plc (LUSE:ZANACO)) removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing values ---
$ws.Range("B2").Value = "'1"
$ws.Range("B2").Style = "Normal"
$ws.Range("I2").Value = 0.002634562967809478
$ws.Range("J2").Value = 0.002634562967809478
$ws.Range("K2").Value = -2.07
$ws.Range("L2").Value = -0.1952830188679245
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("U2").Value = 8.93
$ws.Range("V2").Value = 0.8345794392523365
$ws.Range("X2").Value = 0.1086519618817497
$ws.Range("Z2").Value = 2.070163641201522
$ws.Range("AA2").Value = 0.005453976466415156
$ws.Range("AB2").Value = 0.1027341724849342
$ws.Range("AC2").Value = -0.09728019601851909
$ws.Range("AD2").Value = 3.54
$ws.Range("AE2").Value = 0.8803681627060976
$ws.Range("AF2").Value = 4.420368162706097
$ws.Range("AG2").Value = -4.509631837293902
$ws.Range("AH2").Value = 0.2923452732856594
$ws.Range("AI2").Value = 0.569608048230109
$ws.Range("AJ2").Value = -0.7284917017475958
$ws.Range("AK2").Value = 3.855599423257433
$ws.Range("AN2").Value = 17.35294117647059
$ws.Range("AP2").Value = -22.10603841810737
$ws.Range("W2").Value = -3.631578947368421
$ws.Range("Y2").Value = -3.740230909250171

# remove cells no longer present
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("T2").ClearContents()

# --- Row 3 ---
$ws.Range("B3").Value = "Cavmont Capital Holdings Zambia Plc (LUSE:CCHZ)"
$ws.Range("I3").Value = 0.002634562967809478
$ws.Range("J3").Value = 0.002634562967809478
$ws.Range("K3").Value = -2.07
$ws.Range("L3").Value = -0.1952830188679245
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 8.93
$ws.Range("V3").Value = 0.8345794392523365
$ws.Range("X3").Value = 0.1086519618817497
$ws.Range("Z3").Value = 2.070163641201522
$ws.Range("AA3").Value = 0.005453976466415156
$ws.Range("AB3").Value = 0.1027341724849342
$ws.Range("AC3").Value = -0.09728019601851909
$ws.Range("AD3").Value = 3.54
$ws.Range("AE3").Value = 0.8803681627060976
$ws.Range("AF3").Value = 4.420368162706097
$ws.Range("AG3").Value = -4.509631837293902
$ws.Range("AH3").Value = 0.2923452732856594
$ws.Range("AI3").Value = 0.569608048230109
$ws.Range("AJ3").Value = -0.7284917017475958
$ws.Range("AK3").Value = 3.855599423257433
$ws.Range("AN3").Value = 17.35294117647059
$ws.Range("AP3").Value = -22.10603841810737
$ws.Range("W3").Value = -3.631578947368421
$ws.Range("Y3").Value = -3.740230909250171

# remove cells no longer present
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("T3").ClearContents()

# --- delete row 4 entirely ---
$ws.Rows("4:4").Delete()
